$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$answer = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"

$ws.Range("A26").Value = "Wanneer zijn jullie open?"
$ws.Range("B26").Value = "mailmind.test@zohomail.eu"
$ws.Range("C26").Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Range("D26").Value = "Openingstijden / Locatie"
$ws.Range("E26").Value = $answer
$ws.Range("F26").Value = "2025-06-26 22:26:13"
$ws.Range("G26").Value = "Ja"
$ws.Range("H26").Value = "Nee"
$ws.Range("I26").Value = "Ja"

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A3").Value = "Openingstijden / Locatie"
$dash.Range("B3").Value = 3
$dash.Range("A4").Value = "Productinformatie"
$dash.Range("B4").Value = 2
$dash.Range("A5").Value = "Offerte / Prijsaanvraag"
$dash.Range("B5").Value = 2

# Extend the conditional-formatting ranges to cover the newly added row 26
$ws.Range("D2:D25").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D26"))
$ws.Range("G2:G25").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G26"))
$ws.Range("H2:H25").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H26"))
$ws.Range("I2:I25").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I26"))
